# Microservices 101 / 102 / 103 update script
#
# 1. Slide 3 ("topics"): merge the two runs "Microservices-" + "architecture"
#    into a single run "Microservicesarchitecture" (keeping the err="1"
#    formatting that was on the second run).
# 2. Slide 4 ("microservices-architecture"): rename the two "customer data"
#    textboxes to "customer info" and shrink them (spAutoFit already takes
#    care of width normally, but we pin the exact resulting extent so the
#    two boxes end up the same width).

$p = $ppt.ActivePresentation

# --- Slide 3: "Microservices-" + "architecture" -> "Microservicesarchitecture"
$slide3 = $p.Slides.Item(3)
$contentShape = $slide3.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

# Remove the text of the first run ("Microservices-") entirely. What remains
# is a single run ("architecture") carrying its own (err="1") formatting.
$firstRun = $tr.Characters(1, 14)
$firstRun.Text = ""

# Prepend "Microservices" onto the remaining run so it stays one run with the
# formatting of the (former) second run.
$remaining = $tr.Characters(1, 12)
$remaining.Text = "Microservices" + $remaining.Text

# --- Slide 4: "customer data" textboxes -> "customer info" (+ resize)
$slide4 = $p.Slides.Item(4)

$custBox1 = $slide4.Shapes.Item(22)
$custBox1.TextFrame.TextRange.Text = "customer info"
$custBox1.Width = 1481688 / 12700

$custBox2 = $slide4.Shapes.Item(27)
$custBox2.TextFrame.TextRange.Text = "customer info"
$custBox2.Width = 1481688 / 12700
